$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-332)
# from serial date 45203 (2023-10-04) to 45204 (2023-10-05)
$ws.Range("C2:C332").Value = 45204
